$d = $word.ActiveDocument

$pairs = @(
    @("96×26=2496", "77×26=2002"),
    @("66×39=2574", "20×56=1120"),
    @("31×16=496", "87×26=2262"),
    @("81×94=7614", "27×38=1026"),
    @("77×63=4851", "97×53=5141"),
    @("58×73=4234", "77×81=6237"),
    @("63×65=4095", "28×66=1848"),
    @("41×95=3895", "51×45=2295"),
    @("65×32=2080", "24×12=288"),
    @("31×81=2511", "24×94=2256"),
    @("30×80=2400", "39×57=2223"),
    @("15×91=1365", "20×52=1040"),
    @("76×45=3420", "97×45=4365"),
    @("64×42=2688", "59×60=3540"),
    @("61×51=3111", "78×17=1326"),
    @("22×59=1298", "15×68=1020"),
    @("42×92=3864", "20×75=1500"),
    @("84×55=4620", "15×45=675"),
    @("47×17=799", "25×35=875"),
    @("50×92=4600", "79×32=2528"),
    @("29×82=2378", "36×41=1476"),
    @("53×63=3339", "22×55=1210"),
    @("50×23=1150", "71×67=4757"),
    @("37×97=3589", "31×47=1457"),
    @("17×88=1496", "16×21=336")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
